$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.688.99'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '1.851.55'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '263.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5377'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3202'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07067'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.06'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7762'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07830'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("D13").Value = '1.854.36'
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '89.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.053'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9990'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008036'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9996'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '26.705.20'
$ws.Range("E20").Value = '  +1.04%  '
$ws.Range("D21").Value = '2.074.91'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.657'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.066'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.448'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.229'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.700'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.334'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.89%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08774'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.131'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04884'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7382'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.32%  '
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.869'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.111'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.372'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01756'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4853'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9116'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.45'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.61%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.922'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9996'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.755'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4224'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("E47").Value = '  +1.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.112'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.14'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05852'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9011'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.52%  '
